$d = $word.ActiveDocument

# Mapping of old text -> new text, derived from the unified diff.
$replacements = [ordered]@{
    "2025-08-22 Friday" = "2025-08-23 Saturday"
    "38×88="            = "43×77="
    "68×25="            = "77×99="
    "59×54="            = "48×55="
    "77×61="            = "50×97="
    "83×19="            = "99×30="
    "20×95="            = "54×14="
    "34×60="            = "62×96="
    "63×78="            = "25×94="
    "87×36="            = "46×38="
    "71×40="            = "42×69="
    "27×49="            = "29×55="
    "16×16="            = "36×30="
    "66×84="            = "65×91="
    "59×12="            = "15×58="
    "71×93="            = "39×43="
    "27×51="            = "44×26="
    "59×90="            = "55×46="
    "54×40="            = "23×48="
    "49×72="            = "17×34="
    "94×33="            = "68×37="
    "89×59="            = "80×68="
    "25×35="            = "60×11="
    "69×70="            = "17×99="
    "42×82="            = "71×77="
    "85×94="            = "58×15="
}

foreach ($old in $replacements.Keys) {
    $new = $replacements[$old]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
